# Apply the "en-ru recognizer" reference-table update:
# - reorder/refresh the id -> author-name rows
# - add a "no id yet" marker row ("-") with a freshly generated multi-author
#   candidate list for a person without an id
# - add a worked example row showing a generated id (stored as text, like
#   a user-typed value) reusing an existing id/name pair
# - drop the three rows that no longer apply

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 2-19: ids + matching alternative-name strings ---
$ws.Cells.Item(2, 1).Value2 = 153485
$ws.Cells.Item(2, 2).Value2 = "Kuklina V., Куклина Вера Владимировна"
$ws.Cells.Item(3, 1).Value2 = 484931
$ws.Cells.Item(3, 2).Value2 = "Шустер В.Л., Шустер Владимир Львович"
$ws.Cells.Item(4, 1).Value2 = 534346
$ws.Cells.Item(4, 2).Value2 = "Пунанова С.А., Пунанова Светлана Александровна"
$ws.Cells.Item(5, 1).Value2 = 804891
$ws.Cells.Item(5, 2).Value2 = "Kuzmin V.A., Кузьмин В.А."
$ws.Cells.Item(6, 1).Value2 = 1090961
$ws.Cells.Item(6, 2).Value2 = "Колоколова И.В., Колоколова Ирина Владимировна"
$ws.Cells.Item(7, 1).Value2 = 9036402995
$ws.Cells.Item(7, 2).Value2 = "Краус З.Т., Краус Зоя Тимофеевна"
$ws.Cells.Item(8, 1).Value2 = 1487442777
$ws.Cells.Item(8, 2).Value2 = "Казанин А.Г., Казанин Алексей Геннадьевич"
$ws.Cells.Item(9, 1).Value2 = 1509268305
$ws.Cells.Item(9, 2).Value2 = "Ростовщиков В.Б., Ростовщиков Владимир Борисович"
$ws.Cells.Item(10, 1).Value2 = 2240719343
$ws.Cells.Item(10, 2).Value2 = "Черных С.П., Черных Сергей Петрович"
$ws.Cells.Item(11, 1).Value2 = 4006277815
$ws.Cells.Item(11, 2).Value2 = "Еременко В.Б., Еременко Василий Борисович"
$ws.Cells.Item(12, 1).Value2 = 4247549511
$ws.Cells.Item(12, 2).Value2 = "Кособреева А.А., Кособреева Александра Александровна"
$ws.Cells.Item(13, 1).Value2 = 4786099333
$ws.Cells.Item(13, 2).Value2 = "Замрий А.В., Замрий Анатолий Владимирович"
$ws.Cells.Item(14, 1).Value2 = 5945527727
$ws.Cells.Item(14, 2).Value2 = "Бобов Д.Г., Бобов Дмитрий Геннадиевич, Бобов Дмитрий Геннадьевич"
$ws.Cells.Item(15, 1).Value2 = 8929086524
$ws.Cells.Item(15, 2).Value2 = "Туманова Е.С., Туманова Екатерина Сергеевна"
$ws.Cells.Item(16, 1).Value2 = 9838921473
$ws.Cells.Item(16, 2).Value2 = "Скворцов А.С., Скворцов Андрей Сергеевич"
$ws.Cells.Item(17, 1).Value2 = 6893755920
$ws.Cells.Item(17, 2).Value2 = "Замрий А.В., Замрий Анатолий Владимирович"
$ws.Cells.Item(18, 1).Value2 = 5277594628
$ws.Cells.Item(18, 2).Value2 = "Черных С.П., Черных Сергей Петрович"
$ws.Cells.Item(19, 1).Value2 = 6524091593
$ws.Cells.Item(19, 2).Value2 = "Замрий А.В., Замрий Анатолий Владимирович"

# --- row 20: person with no id yet -> "-" placeholder + generated candidate list ---
$ws.Cells.Item(20, 1).Value2 = "-"
$ws.Cells.Item(20, 2).Value2 = "Бабич Е.А., Гумерова Р.Р., Еремина И.А., Зарипова Л.Ф., Краус З.Т., Краус Зоя Тимофеевна, Лаптев Ярослав Андреевич, Насекин К.К., Пахомов Андрей Львович, Сейнароев М.Р., Суслов А.А., Фицнер Л.К., Чащина-Семенова О.К."

# --- row 21: worked example of a generated id re-entered as typed text ---
# (mirrors how a user would retype an id string into the sheet, so it needs
# to land as literal text "153485" rather than the number 153485; build it
# via a scratch formula cell + paste-values so it comes through as text)
$ws.Cells.Item(200, 1).Formula = "=""153485"""
$ws.Cells.Item(200, 1).Copy()
$ws.Cells.Item(21, 1).PasteSpecial(-4163)
$ws.Cells.Item(200, 1).Clear()
$ws.Cells.Item(21, 2).Value2 = "Kuklina V., Куклина Вера Владимировна"

# --- the table shrank by three rows; drop the old trailing rows ---
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
